# Update countries & provincias Spain
# Refresh the "Pais" COVID-19 table: bump the timestamp, update several
# countries' statistics, and re-rank rows whose totals changed enough to
# move them past a neighboring row (names on those rows swap accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 30 de Marzo de 2020 a las 16:50'

$ws.Range("B4").Value = 144410
$ws.Range("C4").Value = 919
$ws.Range("D4").Value = 4573
$ws.Range("E4").Value = 137237
$ws.Range("F4").Value = 2970
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 2600

$ws.Range("B11").Value = 19522
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 135
$ws.Range("E11").Value = 17972
$ws.Range("F11").Value = 163
$ws.Range("G11").Value = 187
$ws.Range("H11").Value = 1415

$ws.Range("B43").Value = 1156
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 52
$ws.Range("E43").Value = 1063
$ws.Range("F43").Value = 66
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 41

$ws.Range("A48").Value = 'Republica Dominicana'
$ws.Range("B48").Value = 901
$ws.Range("C48").Value = 42
$ws.Range("D48").Value = 4
$ws.Range("E48").Value = 855
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 42

$ws.Range("A49").Value = 'Singapur'
$ws.Range("B49").Value = 879
$ws.Range("C49").Value = 35
$ws.Range("D49").Value = 228
$ws.Range("E49").Value = 648
$ws.Range("F49").Value = 19
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 3

$ws.Range("B73").Value = 354
$ws.Range("C73").Value = 8
$ws.Range("D73").Value = 17
$ws.Range("E73").Value = 329
$ws.Range("F73").Value = 13
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 8

$ws.Range("A135").Value = 'Zambia'
$ws.Range("B135").Value = 35
$ws.Range("C135").Value = 6
$ws.Range("D135").Value = 0
$ws.Range("E135").Value = 35
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

$ws.Range("A136").Value = 'Polinesia Francesa'
$ws.Range("B136").Value = 35
$ws.Range("C136").Value = 5
$ws.Range("D136").Value = 0
$ws.Range("E136").Value = 35
$ws.Range("F136").Value = 2
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

$ws.Range("A138").Value = 'Uganda'

$ws.Range("A139").Value = 'Barbados'

$ws.Range("A164").Value = 'Curazao'
$ws.Range("B164").Value = 11
$ws.Range("C164").Value = 3
$ws.Range("D164").Value = 2
$ws.Range("E164").Value = 8
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

$ws.Range("A165").Value = 'Groenlandia'
$ws.Range("B165").Value = 10
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 2
$ws.Range("E165").Value = 8
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0

$ws.Range("A166").Value = 'Granada'

$ws.Range("A167").Value = 'Suazilandia'
$ws.Range("B167").Value = 9
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 0
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

$ws.Range("A168").Value = 'Santa Lucia'
$ws.Range("B168").Value = 9
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 1
$ws.Range("E168").Value = 8
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

$ws.Range("A169").Value = 'Siria'
$ws.Range("B169").Value = 9
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 8
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 1

$ws.Range("A170").Value = 'Laos'

$ws.Range("A171").Value = 'Seychelles'

$ws.Range("A172").Value = 'Surinam'

$ws.Range("A173").Value = 'Mozambique'

$ws.Range("A174").Value = 'Libia'
$ws.Range("B174").Value = 8
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 8
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

$ws.Range("A175").Value = 'Guinea-Bisau'
$ws.Range("B175").Value = 8
$ws.Range("C175").Value = 6
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 8
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

$ws.Range("A176").Value = 'Guyana'
$ws.Range("B176").Value = 8
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 7
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 1

$ws.Range("A177").Value = 'San Cristobal y Nieves'
$ws.Range("B177").Value = 7
$ws.Range("C177").Value = 5
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 7
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

$ws.Range("A178").Value = 'Antigua y Barbuda'
$ws.Range("B178").Value = 7
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 7
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

$ws.Range("A179").Value = 'Gabon'

$ws.Range("A180").Value = 'Zimbabue'
$ws.Range("B180").Value = 7
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 6
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 1

$ws.Range("A181").Value = 'Angola'
$ws.Range("B181").Value = 7
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 5
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 2

$ws.Range("A184").Value = 'Santa Sede'
$ws.Range("B184").Value = 6
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 6
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

$ws.Range("A185").Value = 'Cabo Verde'
$ws.Range("B185").Value = 6
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 5
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 1

$ws.Range("A186").Value = 'Sudan'
$ws.Range("B186").Value = 6
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 4
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 1
$ws.Range("H186").Value = 2

$ws.Range("A187").Value = 'San Bartolome'

$ws.Range("A188").Value = 'Montserrat'

$ws.Range("A189").Value = 'Fiyi'
$ws.Range("B189").Value = 5
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 5
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

$ws.Range("A190").Value = 'Republica del Chad'
$ws.Range("B190").Value = 5
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 5
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

$ws.Range("A191").Value = 'Nepal'
$ws.Range("B191").Value = 5
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 1
$ws.Range("E191").Value = 4
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

$ws.Range("A192").Value = 'Mauritania'
$ws.Range("B192").Value = 5
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 2
$ws.Range("E192").Value = 3
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = 'Butan'

$ws.Range("A194").Value = 'Islas Turcas y Caicos'
$ws.Range("B194").Value = 4
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 4
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

$ws.Range("A196").Value = 'Gambia'
$ws.Range("B196").Value = 4
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 3
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 1

$ws.Range("A199").Value = 'Republica de Africa Central'
$ws.Range("B199").Value = 3
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 0
$ws.Range("E199").Value = 3
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

$ws.Range("A201").Value = 'Islas Virgenes Britanicas'

$ws.Range("A202").Value = 'Belice'
